# Weekly update: insert the newest week of "Piña" price data (Terminal La
# Palmera de La Serena) at the top of the data block (row 457), pushing the
# previously-existing rows (old 457-466) down by four rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows starting at row 457 (shifts old rows 457:466 -> 461:470)
$ws.Rows("457:460").Insert()

# New week's rows (date 2021-09-09, serial 44448) - one row per "Calidad" grade.
$newRows = @(
    @{ Row=457; A=8; B="Terminal La Palmera de La Serena"; C="Coquimbo"; D=44448; E=4; F="Fruta"; G=100108; H="Tropicales y subtropicales"; I=100108005; J="Piña"; K="Caramelo"; L="Especial"; M=216; N=21000; O=22000; P=21500; Q="$/caja 10 unidades"; R="Ecuador"; S=2150; T=10 },
    @{ Row=458; A=8; B="Terminal La Palmera de La Serena"; C="Coquimbo"; D=44448; E=4; F="Fruta"; G=100108; H="Tropicales y subtropicales"; I=100108005; J="Piña"; K="Caramelo"; L="Primera"; M=216; N=21000; O=22000; P=21500; Q="$/caja 12 unidades"; R="Ecuador"; S=1792; T=12 },
    @{ Row=459; A=8; B="Terminal La Palmera de La Serena"; C="Coquimbo"; D=44448; E=4; F="Fruta"; G=100108; H="Tropicales y subtropicales"; I=100108005; J="Piña"; K="Caramelo"; L="Segunda"; M=216; N=21000; O=22000; P=21500; Q="$/caja 14 unidades"; R="Ecuador"; S=1536; T=14 },
    @{ Row=460; A=8; B="Terminal La Palmera de La Serena"; C="Coquimbo"; D=44448; E=4; F="Fruta"; G=100108; H="Tropicales y subtropicales"; I=100108005; J="Piña"; K="Caramelo"; L="Tercera"; M=216; N=21000; O=22000; P=21500; Q="$/caja 16 unidades"; R="Ecuador"; S=1344; T=16 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
